$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date serial number that was bumped by
# one day (45205 -> 45206) for every data row (rows 2 through 306) as
# part of an automatic update.
$ws.Range("C2:C306").Value = 45206
